$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.747.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.307.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'547.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'131.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.73%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.92%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.307.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.86%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.39%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -4.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.720.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'58.737.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.92%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.13%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.294.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.74%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.21%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.94%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'314.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.77%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'63.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -6.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -5.71%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -5.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.19%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'168.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0726"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.37%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -5.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.76%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.32%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.59%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -5.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -5.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'297.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -7.90%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'141.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.55%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.86%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0502"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.16%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.558"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'18.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.52%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'16.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.35%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.33%  "
$ws.Range("E51").Style = "Normal"
